$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.128.57'
$ws.Range("E2").Value = '  -2.21%  '
$ws.Range("D3").Value = '2.171.45'
$ws.Range("E3").Value = '  -2.17%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.99'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("E6").Value = '  -3.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.90'
$ws.Range("E7").Value = '  -5.13%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.577'
$ws.Range("E9").Value = '  -6.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.06'
$ws.Range("E10").Value = '  -8.25%  '
$ws.Range("E11").Value = '  -3.26%  '
$ws.Range("E12").Value = '  -5.06%  '
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").Value = '2.496.83'
$ws.Range("E14").Value = '  -2.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.90'
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.808'
$ws.Range("E16").Value = '  -4.15%  '
$ws.Range("D17").Value = '2.160.70'
$ws.Range("E17").Value = '  -2.35%  '
$ws.Range("D18").Value = '40.971.75'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("E19").Value = '  -7.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.37'
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("E21").Value = '  -4.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.94'
$ws.Range("E22").Value = '  -5.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '225.55'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  -7.42%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.89'
$ws.Range("E26").Value = '  -5.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.53'
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.77'
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '30.73'
$ws.Range("E32").Value = '  +4.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0770'
$ws.Range("E33").Value = '  -3.51%  '
$ws.Range("E34").Value = '  -8.62%  '
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.104'
$ws.Range("E36").Value = '  -8.43%  '
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '12.40'
$ws.Range("E39").Value = '  -4.31%  '
$ws.Range("E40").Value = '  -3.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.41'
$ws.Range("E41").Value = '  -4.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '60.14'
$ws.Range("E42").Value = '  -7.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.189'
$ws.Range("E43").Value = '  -4.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.31'
$ws.Range("E44").Value = '  -4.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0971'
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '98.16'
$ws.Range("E46").Value = '  -6.21%  '
$ws.Range("E47").Value = '  -2.55%  '
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("E49").Value = '  -7.82%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.63'
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("D51").Value = '2.373.79'
